$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20

# Row 2: Sending=ECs, Ligand=Sema4d, Receptor=Met, Target=ECs
$data[0,0] = "ECs"
$data[0,1] = "Sema4d"
$data[0,2] = "Met"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.8730476666666668
$data[0,7] = 2.619143
$data[0,8] = 0.01740928848427011
$data[0,9] = 0.01740928848427011
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 3.167735333333333
$data[0,13] = 9.503206
$data[0,14] = 0.1182666224938439
$data[0,15] = 0.1182666224938439
$data[0,16] = 2.765583941384222
$data[0,17] = 24.890255472458
$data[0,18] = 0.002058937749055596
$data[0,19] = 0.002058937749055597

# Row 3: Sending=ECs, Ligand=Sema4d, Receptor=Met, Target=FAPs
$data[1,0] = "ECs"
$data[1,1] = "Sema4d"
$data[1,2] = "Met"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.8730476666666668
$data[1,7] = 2.619143
$data[1,8] = 0.01740928848427011
$data[1,9] = 0.01740928848427011
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.9421210000000001
$data[1,13] = 2.826363
$data[1,14] = 0.03517385669126484
$data[1,15] = 0.03517385669126484
$data[1,16] = 0.8225165407676669
$data[1,17] = 7.402648866909001
$data[1,18] = 0.000612351818242604
$data[1,19] = 0.0006123518182426041

# Row 4: Sending=ECs, Ligand=Sema4d, Receptor=Met, Target=M2
$data[2,0] = "ECs"
$data[2,1] = "Sema4d"
$data[2,2] = "Met"
$data[2,3] = "M2"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.8730476666666668
$data[2,7] = 2.619143
$data[2,8] = 0.01740928848427011
$data[2,9] = 0.01740928848427011
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 1.442127
$data[2,13] = 4.326381
$data[2,14] = 0.05384145818700961
$data[2,15] = 0.0538414581870096
$data[2,16] = 1.259045612387
$data[2,17] = 11.331410511483
$data[2,18] = 0.0009373414779914168
$data[2,19] = 0.0009373414779914169

# Row 5: Sending=ECs, Ligand=Sema4d, Receptor=Met, Target=sCs
$data[3,0] = "ECs"
$data[3,1] = "Sema4d"
$data[3,2] = "Met"
$data[3,3] = "sCs"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.8730476666666668
$data[3,7] = 2.619143
$data[3,8] = 0.01740928848427011
$data[3,9] = 0.01740928848427011
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 21.232711
$data[3,13] = 63.69813300000001
$data[3,14] = 0.7927180626278817
$data[3,15] = 0.7927180626278817
$data[3,16] = 18.53716879555767
$data[3,17] = 166.834519160019
$data[3,18] = 0.01380065743898049
$data[3,19] = 0.01380065743898049

# Row 6: Sending=FAPs, Ligand=Sema4d, Receptor=Met, Target=ECs
$data[4,0] = "FAPs"
$data[4,1] = "Sema4d"
$data[4,2] = "Met"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.252512666666667
$data[4,7] = 3.757538
$data[4,8] = 0.02497613266347325
$data[4,9] = 0.02497613266347325
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 3.167735333333333
$data[4,13] = 9.503206
$data[4,14] = 0.1182666224938439
$data[4,15] = 0.1182666224938439
$data[4,16] = 3.967628629647556
$data[4,17] = 35.708657666828
$data[4,18] = 0.002953842853067155
$data[4,19] = 0.002953842853067155

# Row 7: Sending=FAPs, Ligand=Sema4d, Receptor=Met, Target=FAPs
$data[5,0] = "FAPs"
$data[5,1] = "Sema4d"
$data[5,2] = "Met"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.252512666666667
$data[5,7] = 3.757538
$data[5,8] = 0.02497613266347325
$data[5,9] = 0.02497613266347325
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.9421210000000001
$data[5,13] = 2.826363
$data[5,14] = 0.03517385669126484
$data[5,15] = 0.03517385669126484
$data[5,16] = 1.180018486032667
$data[5,17] = 10.620166374294
$data[5,18] = 0.000878506911007027
$data[5,19] = 0.000878506911007027

# Row 8: Sending=FAPs, Ligand=Sema4d, Receptor=Met, Target=M2
$data[6,0] = "FAPs"
$data[6,1] = "Sema4d"
$data[6,2] = "Met"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.252512666666667
$data[6,7] = 3.757538
$data[6,8] = 0.02497613266347325
$data[6,9] = 0.02497613266347325
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 1.442127
$data[6,13] = 4.326381
$data[6,14] = 0.05384145818700961
$data[6,15] = 0.0538414581870096
$data[6,16] = 1.806282334442
$data[6,17] = 16.256541009978
$data[6,18] = 0.0013447514024736
$data[6,19] = 0.0013447514024736

# Row 9: Sending=FAPs, Ligand=Sema4d, Receptor=Met, Target=sCs
$data[7,0] = "FAPs"
$data[7,1] = "Sema4d"
$data[7,2] = "Met"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.252512666666667
$data[7,7] = 3.757538
$data[7,8] = 0.02497613266347325
$data[7,9] = 0.02497613266347325
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 21.232711
$data[7,13] = 63.69813300000001
$data[7,14] = 0.7927180626278817
$data[7,15] = 0.7927180626278817
$data[7,16] = 26.59423947517267
$data[7,17] = 239.348155276554
$data[7,18] = 0.01979903149692547
$data[7,19] = 0.01979903149692547

# Row 10: Sending=M2, Ligand=Sema4d, Receptor=Met, Target=ECs
$data[8,0] = "M2"
$data[8,1] = "Sema4d"
$data[8,2] = "Met"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 45.633473
$data[8,7] = 136.900419
$data[8,8] = 0.9099689814525027
$data[8,9] = 0.9099689814525027
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 3.167735333333333
$data[8,13] = 9.503206
$data[8,14] = 0.1182666224938439
$data[8,15] = 0.1182666224938439
$data[8,16] = 144.5547648048127
$data[8,17] = 1300.992883243314
$data[8,18] = 0.1076189580105508
$data[8,19] = 0.1076189580105508

# Row 11: Sending=M2, Ligand=Sema4d, Receptor=Met, Target=FAPs
$data[9,0] = "M2"
$data[9,1] = "Sema4d"
$data[9,2] = "Met"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 45.633473
$data[9,7] = 136.900419
$data[9,8] = 0.9099689814525027
$data[9,9] = 0.9099689814525027
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.9421210000000001
$data[9,13] = 2.826363
$data[9,14] = 0.03517385669126484
$data[9,15] = 0.03517385669126484
$data[9,16] = 42.99225321623301
$data[9,17] = 386.930278946097
$data[9,18] = 0.03200711854710656
$data[9,19] = 0.03200711854710656

# Row 12: Sending=M2, Ligand=Sema4d, Receptor=Met, Target=M2
$data[10,0] = "M2"
$data[10,1] = "Sema4d"
$data[10,2] = "Met"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 45.633473
$data[10,7] = 136.900419
$data[10,8] = 0.9099689814525027
$data[10,9] = 0.9099689814525027
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 1.442127
$data[10,13] = 4.326381
$data[10,14] = 0.05384145818700961
$data[10,15] = 0.0538414581870096
$data[10,16] = 65.809263517071
$data[10,17] = 592.2833716536389
$data[10,18] = 0.04899405686635065
$data[10,19] = 0.04899405686635065

# Row 13: Sending=M2, Ligand=Sema4d, Receptor=Met, Target=sCs
$data[11,0] = "M2"
$data[11,1] = "Sema4d"
$data[11,2] = "Met"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 45.633473
$data[11,7] = 136.900419
$data[11,8] = 0.9099689814525027
$data[11,9] = 0.9099689814525027
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 21.232711
$data[11,13] = 63.69813300000001
$data[11,14] = 0.7927180626278817
$data[11,15] = 0.7927180626278817
$data[11,16] = 968.9223441353031
$data[11,17] = 8720.301097217727
$data[11,18] = 0.7213488480284947
$data[11,19] = 0.7213488480284947

# Row 14: Sending=sCs, Ligand=Sema4d, Receptor=Met, Target=ECs
$data[12,0] = "sCs"
$data[12,1] = "Sema4d"
$data[12,2] = "Met"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 2.389349666666667
$data[12,7] = 7.168049000000001
$data[12,8] = 0.04764559739975398
$data[12,9] = 0.04764559739975399
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 3.167735333333333
$data[12,13] = 9.503206
$data[12,14] = 0.1182666224938439
$data[12,15] = 0.1182666224938439
$data[12,16] = 7.568827362788223
$data[12,17] = 68.11944626509401
$data[12,18] = 0.005634883881170375
$data[12,19] = 0.005634883881170377

# Row 15: Sending=sCs, Ligand=Sema4d, Receptor=Met, Target=FAPs
$data[13,0] = "sCs"
$data[13,1] = "Sema4d"
$data[13,2] = "Met"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 2.389349666666667
$data[13,7] = 7.168049000000001
$data[13,8] = 0.04764559739975398
$data[13,9] = 0.04764559739975399
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.9421210000000001
$data[13,13] = 2.826363
$data[13,14] = 0.03517385669126484
$data[13,15] = 0.03517385669126484
$data[13,16] = 2.251056497309667
$data[13,17] = 20.259508475787
$data[13,18] = 0.001675879414908647
$data[13,19] = 0.001675879414908648

# Row 16: Sending=sCs, Ligand=Sema4d, Receptor=Met, Target=M2
$data[14,0] = "sCs"
$data[14,1] = "Sema4d"
$data[14,2] = "Met"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 2.389349666666667
$data[14,7] = 7.168049000000001
$data[14,8] = 0.04764559739975398
$data[14,9] = 0.04764559739975399
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 1.442127
$data[14,13] = 4.326381
$data[14,14] = 0.05384145818700961
$data[14,15] = 0.0538414581870096
$data[14,16] = 3.445745666741
$data[14,17] = 31.011711000669
$data[14,18] = 0.002565308440193948
$data[14,19] = 0.002565308440193948

# Row 17: Sending=sCs, Ligand=Sema4d, Receptor=Met, Target=sCs
$data[15,0] = "sCs"
$data[15,1] = "Sema4d"
$data[15,2] = "Met"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 2.389349666666667
$data[15,7] = 7.168049000000001
$data[15,8] = 0.04764559739975398
$data[15,9] = 0.04764559739975399
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 21.232711
$data[15,13] = 63.69813300000001
$data[15,14] = 0.7927180626278817
$data[15,15] = 0.7927180626278817
$data[15,16] = 50.73237095027968
$data[15,17] = 456.5913385525171
$data[15,18] = 0.03776952566348102
$data[15,19] = 0.03776952566348103

$ws.Range("A2:T17").Value = $data
